$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.009355319046315843; C = 0.8453040705977058; D = 3.895540889934523; E = 1.973712463844347; F = 1.98043797059936;  G = 147 }
    3 = @{ B = 0.014819735710887;    C = 0.9348841099345737; D = 3.83779492710257;  E = 1.959029077656217; F = 1.965716495099241; G = 146 }
    4 = @{ B = 0.02337911836098484;  C = 0.9155265504695959; D = 2.856554757174117; E = 1.690134538187454; F = 1.695830641336568; G = 145 }
    5 = @{ B = 0.02567497148425077;  C = 0.9155663396149141; D = 3.371414575219348; E = 1.836141218757247; F = 1.842369964730255; G = 144 }
    6 = @{ B = 0.03261840467407245;  C = 0.9918011744559903; D = 3.939236338740576; E = 1.984750951313685; F = 1.991458256332732; G = 143 }
    7 = @{ B = 0.03961335383754752;  C = 1.002642090472125;  D = 3.553023518865491; E = 1.884946555970617; F = 1.89120118449088;  G = 142 }
    8 = @{ B = 0.05663703026291498;  C = 1.050085343511526;  D = 3.611378132128495; E = 1.900362631743872; F = 1.906290383046426; G = 141 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
